$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet to match the new dataset name
$ws.Name = "dataset"

# New sorted/updated LINESTRING dataset (92 original rows minus 3 removed
# duplicates/segments, plus 11 newly digitized segments, re-sorted A-Z)
$values = @(
    'LINESTRING (106.6525367 10.7902854, 106.6560217 10.7779374)',
    'LINESTRING (106.6560217 10.7779374, 106.6560358 10.7778896, 106.6560711 10.777746, 106.6562253 10.7771342, 106.6563748 10.7766006, 106.6564097 10.7764919, 106.6565488 10.7760081, 106.656688 10.7754926, 106.6567892 10.7751066, 106.656863 10.7748669, 106.6571144 10.7739624, 106.6571379 10.7738452, 106.6572439 10.7734769, 106.6573404 10.773168, 106.6574323 10.7728228, 106.6575329 10.7724473, 106.6576415 10.7720494, 106.6577595 10.7716483, 106.6581076 10.7704288)',
    'LINESTRING (106.6560358 10.7778896, 106.6561196 10.7779083, 106.6562596 10.7780334, 106.6564003 10.7781597, 106.6565919 10.7783381, 106.6568009 10.778577, 106.6568832 10.7786714, 106.6569674 10.7787605, 106.6571359 10.778938, 106.6574642 10.7793056, 106.6576582 10.7795251, 106.6579407 10.7798368, 106.6583734 10.7803017, 106.6587759 10.7807603)',
    'LINESTRING (106.6569368 10.7760085, 106.6565678 10.7759247)',
    'LINESTRING (106.6581076 10.7704288, 106.6586118 10.7685475, 106.659218 10.7663737, 106.6599314 10.7638836)',
    'LINESTRING (106.6588326 10.780762, 106.6636177 10.7764619)',
    'LINESTRING (106.6588326 10.780762, 106.6646163 10.7869224)',
    'LINESTRING (106.6599314 10.7638836, 106.6600421 10.763622, 106.660095 10.7634454, 106.6601554 10.7632761, 106.6603022 10.7627142, 106.660372 10.7623559, 106.660484 10.7620209, 106.6605839 10.7617129, 106.6607267 10.7612478, 106.6608655 10.7607801, 106.6609755 10.7603835, 106.6611076 10.759904, 106.6612283 10.759425, 106.6613416 10.7590449, 106.6614516 10.7586082, 106.6615414 10.7582511, 106.6616762 10.7578183, 106.6618418 10.7571971, 106.6619444 10.7568499, 106.6621523 10.7561786, 106.6623293 10.7555857, 106.6624092 10.755381)',
    'LINESTRING (106.6599314 10.7638836, 106.6670388 10.7677627)',
    'LINESTRING (106.6614079 10.759811, 106.6611475 10.7597572)',
    'LINESTRING (106.6615133 10.7583915, 106.6688963 10.7598514)',
    'LINESTRING (106.6636177 10.7764619, 106.6603428 10.7722665, 106.6584223 10.7706144, 106.6581076 10.7704288)',
    'LINESTRING (106.6636177 10.7764619, 106.6655278 10.7779042, 106.6693926 10.7808404, 106.6720988 10.7829124)',
    'LINESTRING (106.6636177 10.7764619, 106.6670388 10.7677627)',
    'LINESTRING (106.664611 10.786922, 106.6645741 10.7868696, 106.6644641 10.7867524, 106.6642935 10.786562, 106.6641258 10.7863851, 106.6639093 10.7861572, 106.6636521 10.7858871, 106.6633926 10.7856181, 106.6631985 10.7853948, 106.6629437 10.7850993, 106.6626969 10.7847947, 106.6625115 10.7845743, 106.6623428 10.7844169, 106.6621263 10.7841818, 106.6617973 10.7838392, 106.66156 10.783582, 106.6613078 10.7833123, 106.6610426 10.7830323, 106.6608834 10.7828805, 106.6606859 10.7826483, 106.6604492 10.7824711, 106.660361 10.7823996, 106.6600392 10.7820419, 106.6599124 10.7818987, 106.6598122 10.7818657, 106.6596881 10.7817462, 106.6594062 10.7814366, 106.6592094 10.7812215, 106.6589743 10.7809643, 106.6588326 10.780762)',
    'LINESTRING (106.6655039 10.7563308, 106.6655726 10.7560745)',
    'LINESTRING (106.6658816 10.756407, 106.6659393 10.7561586)',
    'LINESTRING (106.6662066 10.7562206, 106.6647449 10.7558855, 106.6637994 10.7556767, 106.6628553 10.7554771, 106.6624092 10.755381)',
    'LINESTRING (106.6662066 10.7562206, 106.6662698 10.7560041, 106.6696064 10.756647, 106.6740375 10.7575061, 106.6776214 10.7588916)',
    'LINESTRING (106.6662066 10.7562206, 106.6674863 10.7566563, 106.6685833 10.7571306, 106.6694094 10.7576444, 106.671351 10.758903, 106.6733037 10.7601928, 106.6745298 10.7610434, 106.6750352 10.7612716, 106.6759269 10.761874, 106.6761978 10.7620703, 106.6764379 10.7621981, 106.6777763 10.7630558, 106.679628 10.7642921, 106.6812266 10.7653303, 106.6818601 10.765624)',
    'LINESTRING (106.6670388 10.7677627, 106.6678592 10.7682218, 106.6707556 10.7697232)',
    'LINESTRING (106.667536 10.7656973, 106.6703043 10.7663246, 106.6717527 10.7666699, 106.6740166 10.7671416)',
    'LINESTRING (106.6678592 10.7682218, 106.6684375 10.7679682, 106.6743695 10.7677038)',
    'LINESTRING (106.6678592 10.7682218, 106.6707556 10.7697232)',
    'LINESTRING (106.668358 10.7610532, 106.668599 10.7611184)',
    'LINESTRING (106.668894 10.7598591, 106.6690388 10.7599757, 106.6692655 10.7600416, 106.6710491 10.7603535, 106.6731748 10.760773, 106.6743992 10.7610554, 106.6745298 10.7610434)',
    'LINESTRING (106.6688963 10.7598514, 106.6665549 10.7567039, 106.6662066 10.7562206)',
    'LINESTRING (106.6688963 10.7598514, 106.6706673 10.7623758, 106.6740166 10.7671416, 106.6743695 10.7677038)',
    'LINESTRING (106.668899 10.7598593, 106.6686611 10.7608182, 106.6684447 10.7618404, 106.6679904 10.7637794, 106.667536 10.7656973, 106.6670388 10.7677627)',
    'LINESTRING (106.6694828 10.7705861, 106.6701215 10.7708758)',
    'LINESTRING (106.6696064 10.756647, 106.6694185 10.7576613, 106.6688963 10.7598514)',
    'LINESTRING (106.6707556 10.7697232, 106.6684126 10.7739895, 106.6674035 10.775555, 106.6655278 10.7779042)',
    'LINESTRING (106.6707556 10.7697232, 106.6730211 10.7709267)',
    'LINESTRING (106.6719337 10.7750366, 106.6718912 10.774609, 106.6718611 10.7743403, 106.6717933 10.7741176, 106.6716746 10.7738423, 106.6716686 10.773781, 106.6716689 10.7737788)',
    'LINESTRING (106.6719337 10.7750366, 106.6749336 10.7748518, 106.6772143 10.774429, 106.6777699 10.7738302, 106.6790228 10.7727415, 106.6820484 10.7700539, 106.684237 10.7681461)',
    'LINESTRING (106.6720988 10.7829124, 106.6646163 10.7869224)',
    'LINESTRING (106.6721033 10.7829304, 106.6712905 10.7833519, 106.6707273 10.7836681, 106.6692722 10.784427, 106.6681202 10.7850382, 106.6674724 10.7853847, 106.6661957 10.786079, 106.664611 10.786922)',
    'LINESTRING (106.6721922 10.7580332, 106.6723805 10.75719)',
    'LINESTRING (106.6730211 10.7709267, 106.6729482 10.7709919, 106.6728762 10.7711058, 106.6724199 10.7719468, 106.6718636 10.773262, 106.6716759 10.7737099, 106.6716686 10.773781)',
    'LINESTRING (106.6740375 10.7575061, 106.6770446 10.7605048, 106.6818601 10.765624)',
    'LINESTRING (106.6743695 10.7677038, 106.6730211 10.7709267)',
    'LINESTRING (106.6743695 10.7677038, 106.6758075 10.7696066, 106.676741 10.7705868, 106.6777173 10.7716514, 106.6790228 10.7727415, 106.6836829 10.7767)',
    'LINESTRING (106.6743695 10.7677038, 106.6764351 10.7621912)',
    'LINESTRING (106.6764351 10.7621912, 106.6770446 10.7605048)',
    'LINESTRING (106.6770446 10.7605048, 106.6776214 10.7588916)',
    'LINESTRING (106.6776214 10.7588916, 106.6833494 10.7611548)',
    'LINESTRING (106.681494 10.777885, 106.6820852 10.7775714, 106.6836091 10.7767448, 106.6851416 10.7759167, 106.685497 10.7757151, 106.6858913 10.7755116, 106.6867953 10.7750341, 106.6874577 10.774673, 106.6882141 10.7742712, 106.6894207 10.7736159, 106.6902405 10.7731118, 106.6910881 10.7726567, 106.6922302 10.7720656, 106.6929755 10.7716407, 106.693131 10.7714678, 106.6932381 10.7713423)',
    'LINESTRING (106.6815678 10.7778402, 106.6720988 10.7829124)',
    'LINESTRING (106.6815678 10.7778402, 106.6811158 10.77735, 106.6798136 10.7752065, 106.6777699 10.7738302, 106.6730211 10.7709267)',
    'LINESTRING (106.6818601 10.765624, 106.6819941 10.7646819, 106.6833494 10.7611548)',
    'LINESTRING (106.6818601 10.765624, 106.6822839 10.7660535, 106.6830349 10.7668229, 106.684237 10.7681461, 106.6845293 10.7683898)',
    'LINESTRING (106.6818601 10.765624, 106.6832867 10.7664882, 106.6844669 10.7672682, 106.6854077 10.7678686, 106.6854704 10.7679114)',
    'LINESTRING (106.6818616 10.7656288, 106.681584 10.7656848, 106.6813795 10.765582, 106.6809182 10.765727, 106.6794966 10.7661624, 106.6786309 10.7664193, 106.6777572 10.7666756, 106.6772556 10.766831, 106.6763121 10.7671222, 106.6757475 10.7673033, 106.6747967 10.7675873, 106.674564 10.7677295, 106.6743695 10.7677038)',
    'LINESTRING (106.6821711 10.7613832, 106.6824021 10.7607771)',
    'LINESTRING (106.6824824 10.7628356, 106.6826757 10.7629067)',
    'LINESTRING (106.6833494 10.7611548, 106.6839404 10.7595122, 106.6844768 10.7581209, 106.6851633 10.756487)',
    'LINESTRING (106.6845293 10.7683898, 106.6848096 10.7686941, 106.6854949 10.7694095, 106.6861346 10.7700841, 106.687005 10.7709958, 106.6872345 10.7712356)',
    'LINESTRING (106.6845293 10.7683898, 106.6854077 10.7678686, 106.6862682 10.7675474, 106.6873643 10.767089, 106.6882854 10.7667126)',
    'LINESTRING (106.6862682 10.7675474, 106.6870243 10.768352, 106.6871303 10.7684298, 106.6874548 10.7688066, 106.6881374 10.7694614, 106.688604 10.7700007)',
    'LINESTRING (106.6872345 10.7712356, 106.6875926 10.7709418, 106.6879024 10.7706836, 106.6879943 10.7705821, 106.6885421 10.7700432, 106.688604 10.7700007)',
    'LINESTRING (106.6872345 10.7712356, 106.6894207 10.7736159)',
    'LINESTRING (106.6882854 10.7667126, 106.6884058 10.7670911, 106.6885185 10.7673889, 106.6887786 10.7678131, 106.6889436 10.7679976, 106.6892762 10.7682123, 106.6899307 10.7684639, 106.6922253 10.7694218, 106.6926638 10.7696036, 106.6932123 10.7698302, 106.6938708 10.7700594, 106.6941944 10.7701595)',
    'LINESTRING (106.6884214 10.7666359, 106.6890115 10.7663764)',
    'LINESTRING (106.688604 10.7700007, 106.6886694 10.770054, 106.688895 10.7702668, 106.689066 10.7704249, 106.6891448 10.7705504, 106.6893064 10.7707385, 106.6894556 10.7709051, 106.6896423 10.7711047, 106.6897801 10.7712536, 106.6899847 10.7714782, 106.6901258 10.7716245, 106.6902579 10.7717668, 106.6904323 10.7719548, 106.6906341 10.7721742, 106.6909171 10.7724759, 106.6910881 10.7726567)',
    'LINESTRING (106.688604 10.7700007, 106.689693 10.7687724, 106.6899546 10.7684746)',
    'LINESTRING (106.6890115 10.7663764, 106.6890363 10.7664304, 106.6890719 10.7665068, 106.6892707 10.7669492)',
    'LINESTRING (106.6890115 10.7663764, 106.689314 10.7662474, 106.6894722 10.7661644, 106.6895366 10.7661644, 106.6896385 10.7660247, 106.6897793 10.765839, 106.6898544 10.7657099, 106.689904 10.7656229, 106.6899646 10.7655183, 106.6900584 10.7653945, 106.6901634 10.7652551)',
    'LINESTRING (106.6892707 10.7669492, 106.689221 10.7670475)',
    'LINESTRING (106.689693 10.7687724, 106.6892762 10.7682123)',
    'LINESTRING (106.6892707 10.7669492, 106.6893115 10.7670064, 106.6893514 10.7670608, 106.6893712 10.7670884, 106.6893846 10.7671085)',
    'LINESTRING (106.689693 10.7687724, 106.6920479 10.7704681, 106.6932381 10.7713423)',
    'LINESTRING (106.6901634 10.7652551, 106.6902135 10.7651824, 106.6902769 10.7650928, 106.6905571 10.7646953)',
    'LINESTRING (106.6901634 10.7652551, 106.690402 10.7654656, 106.6905951 10.7656237, 106.6907574 10.7657515, 106.6908302 10.765806, 106.6909027 10.7658561, 106.6910234 10.7659457, 106.6911132 10.7660181, 106.6911803 10.7660761, 106.6913343 10.7662223)',
    'LINESTRING (106.6905571 10.7646953, 106.6905811 10.7646498, 106.6906022 10.7646097, 106.6906559 10.7645078, 106.6906945 10.7644367, 106.6907172 10.764391, 106.6907232 10.7643324, 106.6907245 10.7642793)',
    'LINESTRING (106.6907245 10.7642793, 106.6906729 10.7642668, 106.6906525 10.764247, 106.6906391 10.7642237, 106.6906337 10.7641983, 106.6906372 10.7641735, 106.6906449 10.764155)',
    'LINESTRING (106.6908026 10.7641975, 106.6913122 10.7643951, 106.6916984 10.7645532, 106.6921437 10.7647245, 106.6923314 10.764569)',
    'LINESTRING (106.6910881 10.7726567, 106.6920886 10.7721692, 106.6931234 10.7716447, 106.6932445 10.7715258, 106.6932381 10.7713423)',
    'LINESTRING (106.6913343 10.7662223, 106.6930596 10.7669878, 106.6940118 10.7674357)',
    'LINESTRING (106.6921437 10.7647245, 106.6922629 10.7647945, 106.6923353 10.7648261, 106.6924399 10.7648367, 106.6925686 10.764842)',
    'LINESTRING (106.6940118 10.7674357, 106.6944743 10.7676544, 106.6947224 10.7677677, 106.6949852 10.7678758, 106.6950837 10.7679157, 106.6951582 10.7679413, 106.6952041 10.767945, 106.6952284 10.767971, 106.6952974 10.7680593, 106.6953075 10.7680731, 106.6953656 10.7680191)',
    'LINESTRING (106.6941944 10.7701595, 106.6940191 10.7704414, 106.6933881 10.7712555, 106.6932381 10.7713423)',
    'LINESTRING (106.6941944 10.7701595, 106.6944271 10.7702555, 106.694688 10.7703509, 106.6951856 10.7705555, 106.6961674 10.7709725, 106.6965876 10.7711707, 106.6967071 10.7712298, 106.6971279 10.7713823, 106.6977198 10.7716192, 106.6978336 10.7716743)',
    'LINESTRING (106.6941944 10.7701595, 106.6946491 10.7693699, 106.6953656 10.7680191)',
    'LINESTRING (106.6944271 10.7702555, 106.6941052 10.7703024)',
    'LINESTRING (106.6947346 10.7712597, 106.6955851 10.7712332, 106.6958573 10.7712161, 106.696324 10.7712227, 106.6965534 10.7712029, 106.6966987 10.7712226)',
    'LINESTRING (106.6947678 10.7712725, 106.6953032 10.7716681, 106.695978 10.772108, 106.6961275 10.7722088, 106.6964341 10.7723983, 106.6971746 10.7728532)',
    'LINESTRING (106.6953656 10.7680191, 106.6943719 10.7668705, 106.6931917 10.7655583, 106.6900696 10.7619536, 106.6869153 10.7583962, 106.6851665 10.7565253)',
    'LINESTRING (106.695978 10.772108, 106.6964474 10.7715257, 106.6966234 10.7713014, 106.6967071 10.7712298)',
    'LINESTRING (106.6959959 10.7720987, 106.6966987 10.7712226)',
    'LINESTRING (106.6964474 10.7715257, 106.6975592 10.7721973)',
    'LINESTRING (106.6965876 10.7711707, 106.6954899 10.7712427, 106.6947678 10.7712725, 106.6938521 10.7713256, 106.6934739 10.771373, 106.6932381 10.7713423)',
    'LINESTRING (106.6971746 10.7728532, 106.6973177 10.7726134, 106.6975592 10.7721973, 106.6978336 10.7716743)',
    'LINESTRING (106.6975382 10.7722206, 106.6964462 10.7715361)',
    'LINESTRING (106.6980208 10.7725168, 106.6975382 10.7722206)',
    'LINESTRING (106.6980208 10.7725168, 106.6976478 10.7731469)',
    'LINESTRING (106.6980208 10.7725168, 106.6983845 10.771998)',
    'LINESTRING (106.6980208 10.7725168, 106.6984918 10.7727937)',
    'LINESTRING (106.6987787 10.7723107, 106.6982533 10.7718855, 106.6978336 10.7716743)',
    'LINESTRING (106.6987815 10.7723253, 106.6984938 10.7727904, 106.6981069 10.7734293)',
    'LINESTRING (106.6990376 10.7740459, 106.6981069 10.7734293, 106.6976461 10.7731452, 106.6971746 10.7728532)'
)

$arr = New-Object 'object[,]' 100,1
for ($i = 0; $i -lt 100; $i++) {
    $arr[$i, 0] = $values[$i]
}

$ws.Range("A1:A100").Value = $arr

# Match the saved selection/scroll position from the edited workbook
$null = $ws.Range("F103").Select()
